# Generate Report for Handoff
# Rotates the per-file rows (the "ffffe7aee18d" file moves into the slot
# previously held by "8cf1f3af", "ffffffa95bafec" moves into the slot
# previously held by "ffffe7aee18d", and "8cf1f3af" - now freshly handed
# off - moves into the last slot with updated status/timestamps/error
# detail), across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md"
$ws.Range("B2").Value = "e2e\ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md"
$ws.Range("G2").Value = "2016-09-02 03:12:30"

$ws.Range("A3").Value = "ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md"
$ws.Range("B3").Value = "e2e\ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md"

$ws.Range("A4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md"
$ws.Range("B4").Value = "e2e\8cf1f3af-dd1c-4115-b698-58fe982d3a89.md"
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-09-02 03:16:00"

# Hyperlinks keep their original r:id -> Address mapping (the rels file is
# unchanged); only which cell / which display text they're attached to
# rotates. Recreate them in original rId order (B2, B3, B4) so the ids
# line up again.
$ws.Range("A1:G10").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md", "", "", "e2e\ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a772eabe95af76628b9a588203a7ba7ea84a8e9/e2e/ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md", "", "", "e2e\ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md", "", "", "e2e\8cf1f3af-dd1c-4115-b698-58fe982d3a89.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("G2").Value = "c91422ca-cea3-4725-bbaf-00ed808af15b.e49ad01d69885edd5ed7cd605418d1d56339b5aa.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-02 03:12:25"
$ws.Range("I2").Value = "c91422ca-cea3-4725-bbaf-00ed808af15b.md"
$ws.Range("J2").Value = "c91422ca-cea3-4725-bbaf-00ed808af15b.e49ad01d69885edd5ed7cd605418d1d56339b5aa.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-02 03:12:43"

$ws.Range("A3").Value = "ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "True"

$ws.Range("A4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.662287748e7b70fc9d976bdd0951905844ddb5bc.zh-cn.xlf"
$ws.Range("H4").Value = "2016-09-02 03:15:55"
$ws.Range("I4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md"
$ws.Range("J4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.662287748e7b70fc9d976bdd0951905844ddb5bc.zh-cn.xlf"
$ws.Range("K4").Value = "2016-09-02 03:15:31"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ae92dcb2a2f8e15cde0969074fa10ccdd69c059/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md."

# Column P ("Error Detail") now needs to hold a long message - widen it.
$ws.Columns.Item(16).ColumnWidth = 39.2

# Hyperlinks: same ref/rId pairing as before, only the display text changes.
$ws.Range("A1:P10").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md", "", "", "ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3416238c3e79311a8e862801052466523a12782c/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md", "", "", "c91422ca-cea3-4725-bbaf-00ed808af15b.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a772eabe95af76628b9a588203a7ba7ea84a8e9/e2e/ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md", "", "", "ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7b26474d107c1827eb26e7fa880466b8903ea51f/e2e/c91422ca-cea3-4725-bbaf-00ed808af15b.md", "", "", "c91422ca-cea3-4725-bbaf-00ed808af15b.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md", "", "", "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7b26474d107c1827eb26e7fa880466b8903ea51f/e2e/c91422ca-cea3-4725-bbaf-00ed808af15b.md", "", "", "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("G2").Value = "c91422ca-cea3-4725-bbaf-00ed808af15b.e49ad01d69885edd5ed7cd605418d1d56339b5aa.de-de.xlf"
$ws.Range("H2").Value = "2016-09-02 03:12:30"
$ws.Range("I2").Value = "c91422ca-cea3-4725-bbaf-00ed808af15b.md"
$ws.Range("J2").Value = "c91422ca-cea3-4725-bbaf-00ed808af15b.e49ad01d69885edd5ed7cd605418d1d56339b5aa.de-de.xlf"
$ws.Range("K2").Value = "2016-09-02 03:12:51"

$ws.Range("A3").Value = "ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "True"

$ws.Range("A4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.662287748e7b70fc9d976bdd0951905844ddb5bc.de-de.xlf"
$ws.Range("H4").Value = "2016-09-02 03:16:00"
$ws.Range("I4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md"
$ws.Range("J4").Value = "8cf1f3af-dd1c-4115-b698-58fe982d3a89.662287748e7b70fc9d976bdd0951905844ddb5bc.de-de.xlf"
$ws.Range("K4").Value = "2016-09-02 03:15:38"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ae92dcb2a2f8e15cde0969074fa10ccdd69c059/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md."

# Column P ("Error Detail") now needs to hold a long message - widen it.
$ws.Columns.Item(16).ColumnWidth = 39.2

# Hyperlinks: same ref/rId pairing as before, only the display text changes.
$ws.Range("A1:P10").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md", "", "", "ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1d368195e133b1bac974b4c6af84377e28488201/e2e/8cf1f3af-dd1c-4115-b698-58fe982d3a89.md", "", "", "c91422ca-cea3-4725-bbaf-00ed808af15b.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a772eabe95af76628b9a588203a7ba7ea84a8e9/e2e/ffffe7aee18d-39a3-4bf0-ad40-1d7a256fef1f.md", "", "", "ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f1eb72526349db65da0faaa65e1f41d61458275c/e2e/c91422ca-cea3-4725-bbaf-00ed808af15b.md", "", "", "c91422ca-cea3-4725-bbaf-00ed808af15b.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99b06ef9243e274a7f894bfb395a8dbe4788b7b8/e2e/ffffffa95bafec-a2c2-4f2b-b56c-7cccf09981d3.md", "", "", "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f1eb72526349db65da0faaa65e1f41d61458275c/e2e/c91422ca-cea3-4725-bbaf-00ed808af15b.md", "", "", "8cf1f3af-dd1c-4115-b698-58fe982d3a89.md")
